# BOM.xlsx update — component value/footprint naming cleanup (resistor units -> "Ohm",
# footprint names standardized) plus a few style tweaks and a changed selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Cell value corrections -------------------------------------------------
$ws.Range("C2").Value  = "0603_C"
$ws.Range("A4").Value  = "10uF"
$ws.Range("A5").Value  = "22Ohm"
$ws.Range("A6").Value  = "10kOhm"
$ws.Range("C7").Value  = "SOT-23-5"
$ws.Range("A8").Value  = "4.7kOhm"
$ws.Range("C9").Value  = "SMD-5032_2P"
$ws.Range("C10").Value = "LED_0805"
$ws.Range("A11").Value = "1.5kOhm"
$ws.Range("A12").Value = "100Ohm"
$ws.Range("A13").Value = "100kOhm"
$ws.Range("A15").Value = "1kOhm"
$ws.Range("C17").Value = "LQFP-48_7X7X05P"
$ws.Range("C18").Value = "LQFP-48_7X7X05P"

# ---- Alignment fixes (center horizontally, matching the rest of the column) -
$a14 = $ws.Range("A14")
$a14.HorizontalAlignment = -4108
$a14.VerticalAlignment = -4108

$a16 = $ws.Range("A16")
$a16.HorizontalAlignment = -4108
$a16.VerticalAlignment = -4108

# ---- New blank, centered cell next to J1 row (C14) --------------------------
$c14 = $ws.Range("C14")
$c14.HorizontalAlignment = -4108
$c14.VerticalAlignment = -4108

# ---- Footprint cells (C17/C18) switch to the small Chinese UI font ----------
foreach ($addr in @("C17", "C18")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Microsoft YaHei"
    $c.Font.Size = 7
    $c.Font.Color = 0
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
}

# ---- C9/D9: highlighted footprint note row (fill + bottom rule + wrap) ------
$c9 = $ws.Range("C9")
$c9.Font.Name = "Microsoft YaHei"
$c9.Font.Size = 7
$c9.Font.Color = 0
$c9.Interior.Pattern = 1
$c9.Interior.Color = 16777215
$c9.Borders.Item(9).LineStyle = 1
$c9.Borders.Item(9).Weight = -4138
$c9.Borders.Item(9).Color = 13750737
$c9.HorizontalAlignment = -4108
$c9.VerticalAlignment = -4108
$c9.WrapText = $true

$d9 = $ws.Range("D9")
$d9.Font.Name = "Microsoft YaHei"
$d9.Font.Size = 7
$d9.Font.Color = 0
$d9.Interior.Pattern = 1
$d9.Interior.Color = 16777215
$d9.Borders.Item(9).LineStyle = 1
$d9.Borders.Item(9).Weight = -4138
$d9.Borders.Item(9).Color = 13750737
$d9.VerticalAlignment = -4160
$d9.WrapText = $true

$ws.Rows.Item(9).RowHeight = 15

# ---- View state: scroll back to top, move the active selection to A12 ------
$w = $excel.ActiveWindow
$w.ScrollRow = 1
$w.ScrollColumn = 1
$ws.Range("A12").Select()
